$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint & Task Tracking")
$summary = $wb.Worksheets.Item("Sprint Summary")

# Helper: write a date-looking value as literal text (not auto-converted
# to a date serial) by forcing Text format during the write, then
# resetting the cell back to the default "Normal" style so no stray
# number-format survives on the saved cell.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# --- Row 13: Task 4.1 ---
$ws.Range("K13").Value = "Done"
Set-TextValue $ws.Range("L13") "2026-02-16"
Set-TextValue $ws.Range("M13") "2026-02-16"
$ws.Range("P13").Value = 'AiProviderService: OpenAI-compatible HTTP via native fetch+AbortController. 10s timeout→deterministic fallback. Response validated. Token usage logged. PromptBuilderService: grading/feedback/doubt prompts (no PII). Cost tiering: grading=cheap, feedback=higher. AiModule @Global. Config: AI_BASE_URL/API_KEY/MODELS/TIMEOUT/MAX_TOKENS/TEMPERATURE.'

# --- Row 14: Task 4.2 ---
$ws.Range("K14").Value = "Done"
Set-TextValue $ws.Range("L14") "2026-02-16"
Set-TextValue $ws.Range("M14") "2026-02-16"
$ws.Range("P14").Value = 'GET /v1/student/activities/:type/:id/feedback?questionId=&level=. FeedbackService: Hint→Approach→Concept→Solution. Auto-advances, monotonic progression enforced. Returns {questionId,level,content,fromAi,nextLevel,maxLevelReached}. Persists feedbackLevel+aiFeedback to response. Static fallback when AI unavailable.'

# --- Row 15: Task 4.3 ---
$ws.Range("K15").Value = "Done"
Set-TextValue $ws.Range("L15") "2026-02-16"
Set-TextValue $ws.Range("M15") "2026-02-16"
$ws.Range("P15").Value = 'GradingService: MCQ/TF/FILL_BLANK→deterministic (case-insensitive, 0/100). SHORT/LONG_ANSWER→AI-assisted with rubric; JSON response parsed. No API key→pending. ActivitiesService delegates to GradingService. AI results stored in response entity. No PII in prompts.'

# --- Row 16: Task 5.1 ---
$ws.Range("K16").Value = "Done"
Set-TextValue $ws.Range("L16") "2026-02-17"
Set-TextValue $ws.Range("M16") "2026-02-17"
$ws.Range("P16").Value = 'GET /v1/student/attendance: AttendanceController + AttendanceQueryDto with period shortcuts (this_month, last_month, this_term) or explicit dates. Summary + calendar. Timezone-safe. Auth + scoped. Validation.'

# --- Row 17: Task 5.2 ---
$ws.Range("K17").Value = "Done"
Set-TextValue $ws.Range("L17") "2026-02-17"
Set-TextValue $ws.Range("M17") "2026-02-17"
$ws.Range("P17").Value = 'DoubtsController: GET list + GET thread + POST message with syllabus context. SyllabusController: GET /v1/student/syllabus/tree. Replies inherit thread syllabus context. AI fallback. Full DTO validation.'

# --- Row 18: Task 5.3 ---
$ws.Range("K18").Value = "Done"
Set-TextValue $ws.Range("L18") "2026-02-17"
Set-TextValue $ws.Range("M18") "2026-02-17"
$ws.Range("P18").Value = 'Enhanced profile with progressOverview (activity breakdown by type with avg scores). GET /v1/student/sync/status with conflict hint (stub - real sync Sprint 7).'

# --- Sprint Summary sheet updates ---
$summary.Range("D2").Value = 8
$summary.Range("F2").Value = 4

$summary.Range("D3").Value = 8
$summary.Range("F3").Value = 3

$summary.Range("D4").Value = 8
$summary.Range("E4").Value = 4
$summary.Range("F4").Value = 4

$summary.Range("F5").Value = 3

$summary.Range("F6").Value = 3
